# feat: add loading feature
# Trim the Yaris_HB suffix lookup table down to the rows that are still
# needed (drop the old rows 2, 4 and 5 - everything below shifts up).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Delete rows from the bottom up so earlier row numbers stay valid.
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(2).Delete()

# Bring the data sheet to the front and select the now-current key column,
# matching the "loading" view the workbook opens to.
$ws1.Activate() | Out-Null
$ws1.Range("A2:A5").Select() | Out-Null
